# Apply the "property boat&car done" edit:
# - Sheet "汽車" (car) row 1 becomes a proper header row (name, capacity, owner,
#   register_date, register_reason, acquire_value, property_category, category,
#   date, legislator_name, legislator_id, source_file, index)
# - Sheet "汽車" row 2 keeps its original data and gains the new metadata columns
#   H2:N2 (property_category=land, category=normal, date=2011-12-28,
#   legislator_name=陳亭妃, legislator_id=1708, source_file=tmp1fdf1, index=27)
# - Sheet "債務" is untouched content-wise (its shared-string indices merely shift
#   on save because of the new strings inserted earlier in the shared table).

$wb = $excel.ActiveWorkbook
# "汽車" (car) is the workbook's tabSelected sheet, i.e. the ActiveSheet.
$carSheet = $wb.ActiveSheet

# Replicate the existing header-row formatting (bold font + border) onto the
# new H1:N1 header cells, and the plain data-row formatting onto H2:N2, by
# copying the format only from the existing analogous cells in that row.
$carSheet.Range("G1").Copy()
$carSheet.Range("H1:N1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$carSheet.Range("G2").Copy()
$carSheet.Range("H2:N2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 1: replace the (erroneously duplicated data) header row with real headers
$carSheet.Range("B1").Value = "name"
$carSheet.Range("C1").Value = "capacity"
$carSheet.Range("D1").Value = "owner"
$carSheet.Range("E1").Value = "register_date"
$carSheet.Range("F1").Value = "register_reason"
$carSheet.Range("G1").Value = "acquire_value"
$carSheet.Range("H1").Value = "property_category"
$carSheet.Range("I1").Value = "category"
$carSheet.Range("J1").Value = "date"
$carSheet.Range("K1").Value = "legislator_name"
$carSheet.Range("L1").Value = "legislator_id"
$carSheet.Range("M1").Value = "source_file"
$carSheet.Range("N1").Value = "index"

# Row 2: append the new metadata columns to the existing data row
$carSheet.Range("H2").Value = "land"
$carSheet.Range("I2").Value = "normal"
$carSheet.Range("K2").Value = "陳亭妃"
$carSheet.Range("L2").Value = 1708
$carSheet.Range("M2").Value = "tmp1fdf1"
$carSheet.Range("N2").Value = 27

# "date" needs to land as plain text "2011-12-28" (not an auto-converted date
# serial) — build it via a text formula, then paste-values it back over itself
# so it collapses to a literal shared string with no leftover formula/number
# formatting.
$carSheet.Range("J2").Formula = '="2011-12-28"'
$carSheet.Range("J2").Copy()
$carSheet.Range("J2").PasteSpecial(-4163)
$excel.CutCopyMode = $false
